$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077005770509684
$ws.Range("D2").Value = 1.07910128515306
$ws.Range("E2").Value = 1.080506282421035
$ws.Range("F2").Value = 1.090788279320382
$ws.Range("I2").Value = 1.057794804003207
$ws.Range("J2").Value = 1.081902390064228
$ws.Range("K2").Value = 1.081778218676892
$ws.Range("L2").Value = 1.083179538724866
$ws.Range("M2").Value = 1.093434936186371
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.078262152919711
$ws.Range("D3").Value = 1.080123569077601
$ws.Range("E3").Value = 1.081619605605099
$ws.Range("F3").Value = 1.091963856306752
$ws.Range("I3").Value = 1.058189823679921
$ws.Range("J3").Value = 1.082817381373135
$ws.Range("K3").Value = 1.08261801192799
$ws.Range("L3").Value = 1.084110416354603
$ws.Range("M3").Value = 1.094429845571146
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.079074988713257
$ws.Range("D4").Value = 1.080784886679979
$ws.Range("E4").Value = 1.082340114499081
$ws.Range("F4").Value = 1.092724739205675
$ws.Range("I4").Value = 1.058444131242422
$ws.Range("J4").Value = 1.083408752860225
$ws.Range("K4").Value = 1.083160632370091
$ws.Range("L4").Value = 1.084712274667444
$ws.Range("M4").Value = 1.095073238099369
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.079416676290542
$ws.Range("D5").Value = 1.081062865215784
$ws.Range("E5").Value = 1.082643044904913
$ws.Range("F5").Value = 1.093044664867026
$ws.Range("I5").Value = 1.05855073227108
$ws.Range("J5").Value = 1.083657201620815
$ws.Range("K5").Value = 1.083388563778727
$ws.Range("L5").Value = 1.084965181618334
$ws.Range("M5").Value = 1.095343630351738
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.079474045518927
$ws.Range("D6").Value = 1.081109536753049
$ws.Range("E6").Value = 1.082693909954061
$ws.Range("F6").Value = 1.093098384788381
$ws.Range("I6").Value = 1.058568612897438
$ws.Range("J6").Value = 1.083698907688431
$ws.Range("K6").Value = 1.083426823573182
$ws.Range("L6").Value = 1.085007639137121
$ws.Range("M6").Value = 1.095389025156865
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.07907955446926
$ws.Range("D7").Value = 1.08078860119635
$ws.Range("E7").Value = 1.082344162156766
$ws.Range("F7").Value = 1.092729013868148
$ws.Range("I7").Value = 1.058445556867494
$ws.Range("J7").Value = 1.083412073286952
$ws.Range("K7").Value = 1.083163678731248
$ws.Range("L7").Value = 1.084715654470242
$ws.Range("M7").Value = 1.095076851446015
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.077430398438902
$ws.Range("D8").Value = 1.079446806240234
$ws.Range("E8").Value = 1.08088251218453
$ws.Range("F8").Value = 1.091185529127079
$ws.Range("I8").Value = 1.057928571539647
$ws.Range("J8").Value = 1.082211759032395
$ws.Range("K8").Value = 1.082062193231961
$ws.Range("L8").Value = 1.083494233562162
$ws.Range("M8").Value = 1.093771249888957
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.074523309072562
$ws.Range("D9").Value = 1.077081059230147
$ws.Range("E9").Value = 1.078307724692237
$ws.Range("F9").Value = 1.088467236221836
$ws.Range("I9").Value = 1.057007622998563
$ws.Range("J9").Value = 1.080091322565268
$ws.Range("K9").Value = 1.080115202330393
$ws.Range("L9").Value = 1.081338196811559
$ws.Range("M9").Value = 1.091467650568301
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072584407486494
$ws.Range("D10").Value = 1.075502936145591
$ws.Range("E10").Value = 1.07659168394266
$ws.Range("F10").Value = 1.086655988928857
$ws.Range("I10").Value = 1.056386929655789
$ws.Range("J10").Value = 1.078674040325063
$ws.Range("K10").Value = 1.078813088750745
$ws.Range("L10").Value = 1.07989825794154
$ws.Range("M10").Value = 1.089929854420749
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.071744613369308
$ws.Range("D11").Value = 1.074819349437819
$ws.Range("E11").Value = 1.075848717499337
$ws.Range("F11").Value = 1.085871905747681
$ws.Range("I11").Value = 1.056116558694187
$ws.Range("J11").Value = 1.078059456532415
$ws.Range("K11").Value = 1.078248267235919
$ws.Range("L11").Value = 1.079274120380863
$ws.Range("M11").Value = 1.089263465691169
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.071432638285627
$ws.Range("D12").Value = 1.074565396125409
$ws.Range("E12").Value = 1.075572758744414
$ws.Range("F12").Value = 1.085580690570619
$ws.Range("I12").Value = 1.056015888729517
$ws.Range("J12").Value = 1.077831037263785
$ws.Range("K12").Value = 1.078038316106227
$ws.Range("L12").Value = 1.079042191103309
$ws.Range("M12").Value = 1.089015860665529
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.071499559767911
$ws.Range("D13").Value = 1.074619871749453
$ws.Range("E13").Value = 1.075631952312979
$ws.Range("F13").Value = 1.085643155963777
$ws.Range("I13").Value = 1.056037493748902
$ws.Range("J13").Value = 1.07788004012872
$ws.Range("K13").Value = 1.078083358203588
$ws.Range("L13").Value = 1.079091945124541
$ws.Range("M13").Value = 1.08906897636459
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.071718826204771
$ws.Range("D14").Value = 1.074798358376638
$ws.Range("E14").Value = 1.0758259064254
$ws.Range("F14").Value = 1.085847833244002
$ws.Range("I14").Value = 1.056108242227773
$ws.Range("J14").Value = 1.078040578097425
$ws.Range("K14").Value = 1.078230915703853
$ws.Range("L14").Value = 1.079254951017572
$ws.Range("M14").Value = 1.089243000194391
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071853918360469
$ws.Range("D15").Value = 1.074908324699879
$ws.Range("E15").Value = 1.075945409507723
$ws.Range("F15").Value = 1.085973945362667
$ws.Range("I15").Value = 1.056151800590155
$ws.Range("J15").Value = 1.078139472837962
$ws.Range("K15").Value = 1.078321810663929
$ws.Range("L15").Value = 1.079355371457558
$ws.Range("M15").Value = 1.089350211560925
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.072640136645866
$ws.Range("D16").Value = 1.075548298181665
$ws.Range("E16").Value = 1.076640993954701
$ws.Range("F16").Value = 1.086708029996294
$ws.Range("I16").Value = 1.05640483937007
$ws.Range("J16").Value = 1.078714809314753
$ws.Range("K16").Value = 1.07885055295625
$ws.Range("L16").Value = 1.079939666401958
$ws.Range("M16").Value = 1.089974069510299
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.073133245363372
$ws.Range("D17").Value = 1.075949669013211
$ws.Range("E17").Value = 1.077077338852384
$ws.Range("F17").Value = 1.08716855396798
$ws.Range("I17").Value = 1.056563133221975
$ws.Range("J17").Value = 1.079075462959167
$ws.Range("K17").Value = 1.079181950935578
$ws.Range("L17").Value = 1.080306008267179
$ws.Range("M17").Value = 1.090365260735187
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073420844648948
$ws.Range("D18").Value = 1.076183757897346
$ws.Range("E18").Value = 1.077331860288495
$ws.Range("F18").Value = 1.087437189107982
$ws.Range("I18").Value = 1.056655308382048
$ws.Range("J18").Value = 1.079285740245737
$ws.Range("K18").Value = 1.079375153596187
$ws.Range("L18").Value = 1.080519628156726
$ws.Range("M18").Value = 1.090593386580333
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073518904801964
$ws.Range("D19").Value = 1.076263572139
$ws.Range("E19").Value = 1.077418647066124
$ws.Range("F19").Value = 1.087528790071242
$ws.Range("I19").Value = 1.05668671146783
$ws.Range("J19").Value = 1.079357424846791
$ws.Range("K19").Value = 1.079441014438919
$ws.Range("L19").Value = 1.080592456695122
$ws.Range("M19").Value = 1.090671163268689
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.073080341831044
$ws.Range("D20").Value = 1.075906608207654
$ws.Range("E20").Value = 1.077030522272197
$ws.Range("F20").Value = 1.087119142120987
$ws.Range("I20").Value = 1.056546165824635
$ws.Range("J20").Value = 1.079036777131941
$ws.Range("K20").Value = 1.079146405027211
$ws.Range("L20").Value = 1.080266709585771
$ws.Range("M20").Value = 1.090323294738373
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.071654258778073
$ws.Range("D21").Value = 1.07474579959487
$ws.Range("E21").Value = 1.075768791452849
$ws.Range("F21").Value = 1.085787560120446
$ws.Range("I21").Value = 1.056087415254854
$ws.Range("J21").Value = 1.077993307419712
$ws.Range("K21").Value = 1.078187467884624
$ws.Range("L21").Value = 1.079206952528778
$ws.Range("M21").Value = 1.089191756695458
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.070757401728997
$ws.Range("D22").Value = 1.074015728144503
$ws.Range("E22").Value = 1.074975560235842
$ws.Range("F22").Value = 1.084950504165097
$ws.Range("I22").Value = 1.055797579169669
$ws.Range("J22").Value = 1.077336451942168
$ws.Range("K22").Value = 1.077583670158446
$ws.Range("L22").Value = 1.078540080434742
$ws.Range("M22").Value = 1.088479858414581
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.071232864521776
$ws.Range("D23").Value = 1.074402774617868
$ws.Range("E23").Value = 1.075396061084094
$ws.Range("F23").Value = 1.085394228512843
$ws.Range("I23").Value = 1.05595135985575
$ws.Range("J23").Value = 1.077684738365738
$ws.Range("K23").Value = 1.077903838239687
$ws.Range("L23").Value = 1.07889365559487
$ws.Range("M23").Value = 1.088857292841863
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.073104246730218
$ws.Range("D24").Value = 1.075926065606205
$ws.Range("E24").Value = 1.077051676642378
$ws.Range("F24").Value = 1.08714146914821
$ws.Range("I24").Value = 1.056553833139997
$ws.Range("J24").Value = 1.079054257857476
$ws.Range("K24").Value = 1.079162466990378
$ws.Range("L24").Value = 1.080284467157241
$ws.Range("M24").Value = 1.09034225751858
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075275000564557
$ws.Range("D25").Value = 1.0776928258892
$ws.Range("E25").Value = 1.078973278376621
$ws.Range("F25").Value = 1.089169807356515
$ws.Range("I25").Value = 1.05724689350608
$ws.Range("J25").Value = 1.080640145467534
$ws.Range("K25").Value = 1.080619267594971
$ws.Range("L25").Value = 1.081896033891164
$ws.Range("M25").Value = 1.092063544631949
